$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos price/volume table. Price cells that would otherwise
# be auto-parsed as numbers by Excel (plain "123.45"-style strings) are
# forced back to text with NumberFormat "@" before the value is written so
# the exact original text formatting (e.g. "352.43") is preserved.
$ws.Range("D2").Value = '51.813.19'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '2.931.25'
$ws.Range("E3").Value = '  +3.69%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '352.43'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.93'
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.559'
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.28'
$ws.Range("E10").Value = '  -2.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0892'
$ws.Range("E11").Value = '  +5.13%  '
$ws.Range("E12").Value = '  +1.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.99'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.397.45'
$ws.Range("E14").Value = '  +3.65%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.74'
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").Value = '2.937.48'
$ws.Range("E16").Value = '  +3.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.982'
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").Value = '51.921.77'
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("E19").Value = '  -4.51%  '
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.22'
$ws.Range("E21").Value = '  +6.67%  '
$ws.Range("D22").Value = '0.0₃0986'
$ws.Range("E22").Value = '  +1.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.27'
$ws.Range("E23").Value = '  +1.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.26'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.79'
$ws.Range("E25").Value = '  +1.19%  '
$ws.Range("E26").Value = '  +11.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.87'
$ws.Range("E27").Value = '  +2.50%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.23'
$ws.Range("E29").Value = '  +14.96%  '
$ws.Range("E30").Value = '  +16.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.54'
$ws.Range("E31").Value = '  -0.50%  '
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.26'
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.25'
$ws.Range("E33").Value = '  +10.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '37.08'
$ws.Range("E34").Value = '  -4.39%  '
$ws.Range("E35").Value = '  +0.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0453'
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.33'
$ws.Range("E38").Value = '  +3.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.75'
$ws.Range("E39").Value = '  -1.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.04'
$ws.Range("E40").Value = '  +1.66%  '
$ws.Range("E41").Value = '  +6.63%  '
$ws.Range("E42").Value = '  +0.97%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.10'
$ws.Range("E43").Value = '  +3.54%  '
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.50'
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Value = '2.170.51'
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '114.13'
$ws.Range("E48").Value = '  -6.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.247'
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0348'
$ws.Range("E50").Value = '  +11.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.937'
$ws.Range("E51").Value = '  -1.50%  '
